$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- New output-file names added to column D (rows 3 and 4) ---
$ws.Range("D3").Value = "ALL_R_01_02_CompareToSource.csv"
$ws.Range("D4").Value = "ALL_R_01_03_STUDYPOP.csv"

# --- Widen column C to fit the new, longer description text ---
$ws.Columns.Item(3).ColumnWidth = 107.5

# --- Move the selection from C12 to C10 ---
$ws.Range("C10").Select()
